$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 46789
$ws.Range("E4").Value = 126189040

$ws.Range("C8").Value = 181369
$ws.Range("E8").Value = 653049910

$ws.Range("C115").Value = 17560
$ws.Range("E115").Value = 38644994

$ws.Range("C152").Value = 126053
$ws.Range("E152").Value = 716083919

$ws.Range("C164").Value = 50585
$ws.Range("E164").Value = 168943069

$ws.Range("C168").Value = 285114
$ws.Range("E168").Value = 1213524897

$ws.Range("C169").Value = 562672
$ws.Range("E169").Value = 1286084501

$ws.Range("C170").Value = 367559
$ws.Range("E170").Value = 2848069929

$ws.Range("C171").Value = 115224
$ws.Range("E171").Value = 448839617

$ws.Range("C174").Value = 357363
$ws.Range("E174").Value = 1019961276

$ws.Range("C175").Value = 125687
$ws.Range("E175").Value = 815654401

$ws.Range("C179").Value = 235796
$ws.Range("E179").Value = 813626581

$ws.Range("C180").Value = 141529
$ws.Range("E180").Value = 341235629

$ws.Range("C220").Value = 4714
$ws.Range("E220").Value = 11992288

$ws.Range("C237").Value = 58310
$ws.Range("E237").Value = 172543972

$ws.Range("C239").Value = 84901
$ws.Range("E239").Value = 500391684

$ws.Range("C255").Value = 141372
$ws.Range("E255").Value = 414550820

$ws.Range("C303").Value = 40039
$ws.Range("E303").Value = 131830133
